$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: week number 44 -> 43
$ws.Range("B1").Value = 43

# A2: date 45957 (27/10/2025, lundi) -> 45951 (21/10/2025, mardi)
$ws.Range("A2").Value = 45951

# B2: day name "lundi" -> "mardi"
$ws.Range("B2").Value = "mardi"

# A3: room/cloud label gains a second course code
$ws.Range("A3").Value = "Cloud - C (KRTA9AA3/KUPT9BB1)"

# B3: course code gains a second course code
$ws.Range("B3").Value = "KRTA9AA3/KUPT9BB1"

# F3: was empty, now holds a room code
$ws.Range("F3").Value = "U3-106"
